$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date column header
$ws.Range("D1").Value = "19-04-2025"

# Mark attendance "P" for the new date in the same rows that have a "P" in column C
$ws.Range("D27").Value = "P"
$ws.Range("D28").Value = "P"
$ws.Range("D43").Value = "P"
$ws.Range("D44").Value = "P"
$ws.Range("D68").Value = "P"
